$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pressure Sensor")

# --- New small calc block near the top (rows 16-17) ---
$ws.Range("R16").Formula = "=1/(1000000)"
$ws.Range("Q17").Formula = "=500000"
$ws.Range("R17").Formula = "=R16*Q17"

# --- First duty-cycle calc block (rows 33-35) ---
$ws.Range("R33").Formula = "=1645-1635"
$ws.Range("T33").Formula = "=R33*R34"

$ws.Range("R34").Formula = "=150-1"
$ws.Range("T34").Formula = "=T33/R35"

$ws.Range("R35").Formula = "=14745-1638"
$ws.Range("T35").Formula = "=T34+1"
$ws.Range("U35").Value = "Above Atmospheric"

# --- Second duty-cycle calc block (rows 37-39) ---
$ws.Range("R37").Formula = "=1673-1635"
$ws.Range("T37").Formula = "=R37*R38"

$ws.Range("R38").Formula = "=150-1"
$ws.Range("T38").Formula = "=T37/R39"

$ws.Range("R39").Formula = "=14745-1638"
$ws.Range("T39").Formula = "=T38+1"
$ws.Range("U39").Value = "Above Atmospheric"

# --- Update the sheet view to match the commit (scrolled/selection state) ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("R15").Select()
